# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G) replacing the old Strike# counts for rows 2-17
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 3
    7  = 4
    8  = 1
    9  = 1
    10 = 3
    11 = 4
    12 = 4
    13 = 5
    14 = 2
    15 = 3
    16 = 2
    17 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
